$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting (style) of the other
# header cells (bold, bordered, centered) from G1 so it reuses the same
# cell style instead of creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add Save values for the two data rows (plain numbers, no special style)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
